$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035492529207973
$ws.Cells.Item(2, 4).Value = 1.037378459628777
$ws.Cells.Item(2, 5).Value = 1.043347320323245
$ws.Cells.Item(2, 6).Value = 1.051433521651894
$ws.Cells.Item(2, 9).Value = 1.036698782795838
$ws.Cells.Item(2, 10).Value = 1.040606006761255
$ws.Cells.Item(2, 11).Value = 1.040169514315052
$ws.Cells.Item(2, 12).Value = 1.046121442126428
$ws.Cells.Item(2, 13).Value = 1.054185036440427

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.036345533137149
$ws.Cells.Item(3, 4).Value = 1.037993403367822
$ws.Cells.Item(3, 5).Value = 1.044140073716059
$ws.Cells.Item(3, 6).Value = 1.052386945588968
$ws.Cells.Item(3, 9).Value = 1.036872791389277
$ws.Cells.Item(3, 10).Value = 1.041103068719073
$ws.Cells.Item(3, 11).Value = 1.040594769499334
$ws.Cells.Item(3, 12).Value = 1.046725259811469
$ws.Cells.Item(3, 13).Value = 1.054950743636264

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.036898069220728
$ws.Cells.Item(4, 4).Value = 1.03839175007892
$ws.Cells.Item(4, 5).Value = 1.044653998068479
$ws.Cells.Item(4, 6).Value = 1.053005191315282
$ws.Cells.Item(4, 9).Value = 1.036984456721124
$ws.Cells.Item(4, 10).Value = 1.041424628648928
$ws.Cells.Item(4, 11).Value = 1.040869681295405
$ws.Cells.Item(4, 12).Value = 1.0471162729038
$ws.Cells.Item(4, 13).Value = 1.055446898853028

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037130494042757
$ws.Cells.Item(5, 4).Value = 1.038559317984575
$ws.Cells.Item(5, 5).Value = 1.04487027995685
$ws.Cells.Item(5, 6).Value = 1.053265415280732
$ws.Cells.Item(5, 9).Value = 1.037031177565089
$ws.Cells.Item(5, 10).Value = 1.041559793909641
$ws.Cells.Item(5, 11).Value = 1.040985191581794
$ws.Cells.Item(5, 12).Value = 1.047280725910375
$ws.Cells.Item(5, 13).Value = 1.055655646497785

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.037169527279693
$ws.Cells.Item(6, 4).Value = 1.038587459361664
$ws.Cells.Item(6, 5).Value = 1.044906607963213
$ws.Cells.Item(6, 6).Value = 1.053309126360177
$ws.Cells.Item(6, 9).Value = 1.037039009089341
$ws.Cells.Item(6, 10).Value = 1.041582487629904
$ws.Cells.Item(6, 11).Value = 1.041004582560401
$ws.Cells.Item(6, 12).Value = 1.047308342423833
$ws.Cells.Item(6, 13).Value = 1.055690705717544

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.03690117434966
$ws.Cells.Item(7, 4).Value = 1.038393988727059
$ws.Cells.Item(7, 5).Value = 1.044656887143564
$ws.Cells.Item(7, 6).Value = 1.053008667211149
$ws.Cells.Item(7, 9).Value = 1.036985081885402
$ws.Cells.Item(7, 10).Value = 1.041426434808087
$ws.Cells.Item(7, 11).Value = 1.040871224997073
$ws.Cells.Item(7, 12).Value = 1.04711847005524
$ws.Cells.Item(7, 13).Value = 1.055449687506834

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.03578068373358
$ws.Cells.Item(8, 4).Value = 1.037586191259892
$ws.Cells.Item(8, 5).Value = 1.043615035262156
$ws.Cells.Item(8, 6).Value = 1.051755462050887
$ws.Cells.Item(8, 9).Value = 1.036757781957621
$ws.Cells.Item(8, 10).Value = 1.04077400544365
$ws.Cells.Item(8, 11).Value = 1.040313283754073
$ws.Cells.Item(8, 12).Value = 1.046325441652573
$ws.Cells.Item(8, 13).Value = 1.054443666671308

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.033810787578096
$ws.Cells.Item(9, 4).Value = 1.0361661694757
$ws.Cells.Item(9, 5).Value = 1.041786586476362
$ws.Cells.Item(9, 6).Value = 1.049557316380466
$ws.Cells.Item(9, 9).Value = 1.036350157612704
$ws.Cells.Item(9, 10).Value = 1.039623840282886
$ws.Cells.Item(9, 11).Value = 1.039328205282475
$ws.Cells.Item(9, 12).Value = 1.044930402497078
$ws.Cells.Item(9, 13).Value = 1.052676293576075

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.032500673915797
$ws.Cells.Item(10, 4).Value = 1.03522188756141
$ws.Cells.Item(10, 5).Value = 1.040572716346617
$ws.Cells.Item(10, 6).Value = 1.048098822853909
$ws.Cells.Item(10, 9).Value = 1.036073678063231
$ws.Cells.Item(10, 10).Value = 1.038856798115761
$ws.Cells.Item(10, 11).Value = 1.038670268545391
$ws.Cells.Item(10, 12).Value = 1.044002056308443
$ws.Cells.Item(10, 13).Value = 1.051501743492427

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.031934146405743
$ws.Cells.Item(11, 4).Value = 1.034813593157688
$ws.Cells.Item(11, 5).Value = 1.040048327328079
$ws.Cells.Item(11, 6).Value = 1.047468946416308
$ws.Cells.Item(11, 9).Value = 1.035952845644265
$ws.Cells.Item(11, 10).Value = 1.038524612513417
$ws.Cells.Item(11, 11).Value = 1.038385100500157
$ws.Cells.Item(11, 12).Value = 1.043600487656663
$ws.Cells.Item(11, 13).Value = 1.050994046991437

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.031723828607579
$ws.Cells.Item(12, 4).Value = 1.034662024161505
$ws.Cells.Item(12, 5).Value = 1.039853731642988
$ws.Cells.Item(12, 6).Value = 1.047235233566522
$ws.Cells.Item(12, 9).Value = 1.035907796237382
$ws.Cells.Item(12, 10).Value = 1.038401217491344
$ws.Cells.Item(12, 11).Value = 1.038279135934239
$ws.Cells.Item(12, 12).Value = 1.043451390352898
$ws.Cells.Item(12, 13).Value = 1.050805601392407

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.0317689372273
$ws.Cells.Item(13, 4).Value = 1.034694532138193
$ws.Cells.Item(13, 5).Value = 1.039895464644746
$ws.Cells.Item(13, 6).Value = 1.047285354365573
$ws.Cells.Item(13, 9).Value = 1.035917467028809
$ws.Cells.Item(13, 10).Value = 1.038427686417167
$ws.Cells.Item(13, 11).Value = 1.03830186750524
$ws.Cells.Item(13, 12).Value = 1.043483369349959
$ws.Cells.Item(13, 13).Value = 1.050846017456994

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.031916759101055
$ws.Cells.Item(14, 4).Value = 1.034801062577
$ws.Cells.Item(14, 5).Value = 1.040032238188849
$ws.Cells.Item(14, 6).Value = 1.047449622502426
$ws.Cells.Item(14, 9).Value = 1.035949125246541
$ws.Cells.Item(14, 10).Value = 1.038514412767995
$ws.Cells.Item(14, 11).Value = 1.0383763422546
$ws.Cells.Item(14, 12).Value = 1.043588161930384
$ws.Cells.Item(14, 13).Value = 1.050978467246528

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032007852418607
$ws.Cells.Item(15, 4).Value = 1.034866711436509
$ws.Cells.Item(15, 5).Value = 1.040116533548475
$ws.Cells.Item(15, 6).Value = 1.047550866876499
$ws.Cells.Item(15, 9).Value = 1.035968608827166
$ws.Cells.Item(15, 10).Value = 1.038567846910004
$ws.Cells.Item(15, 11).Value = 1.038422223279854
$ws.Cells.Item(15, 12).Value = 1.043652736509828
$ws.Cells.Item(15, 13).Value = 1.051060091933202

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.032538288629306
$ws.Cells.Item(16, 4).Value = 1.035248997205655
$ws.Cells.Item(16, 5).Value = 1.040607544271398
$ws.Cells.Item(16, 6).Value = 1.048140660823016
$ws.Cells.Item(16, 9).Value = 1.036081673876986
$ws.Cells.Item(16, 10).Value = 1.038878843202602
$ws.Cells.Item(16, 11).Value = 1.038689188471457
$ws.Cells.Item(16, 12).Value = 1.044028715881339
$ws.Cells.Item(16, 13).Value = 1.051535456551647

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.032871221882252
$ws.Cells.Item(17, 4).Value = 1.035488953030921
$ws.Cells.Item(17, 5).Value = 1.040915871455385
$ws.Cells.Item(17, 6).Value = 1.048511068901718
$ws.Cells.Item(17, 9).Value = 1.036152298544616
$ws.Cells.Item(17, 10).Value = 1.039073910081682
$ws.Cells.Item(17, 11).Value = 1.038856575324956
$ws.Cells.Item(17, 12).Value = 1.044264668910469
$ws.Cells.Item(17, 13).Value = 1.051833879980042

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033065489477312
$ws.Cells.Item(18, 4).Value = 1.035628971620669
$ws.Cells.Item(18, 5).Value = 1.041095831540077
$ws.Cells.Item(18, 6).Value = 1.048727281799223
$ws.Cells.Item(18, 9).Value = 1.036193385098729
$ws.Cells.Item(18, 10).Value = 1.039187684263061
$ws.Cells.Item(18, 11).Value = 1.03895418243012
$ws.Cells.Item(18, 12).Value = 1.044402335920105
$ws.Cells.Item(18, 13).Value = 1.052008031222919

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.033131742087444
$ws.Cells.Item(19, 4).Value = 1.035676723858841
$ws.Cells.Item(19, 5).Value = 1.041157213249613
$ws.Cells.Item(19, 6).Value = 1.048801031917027
$ws.Cells.Item(19, 9).Value = 1.036207376278246
$ws.Cells.Item(19, 10).Value = 1.039226477429063
$ws.Cells.Item(19, 11).Value = 1.038987459346676
$ws.Cells.Item(19, 12).Value = 1.044449283487929
$ws.Cells.Item(19, 13).Value = 1.052067426820041

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.032835493706391
$ws.Cells.Item(20, 4).Value = 1.035463202194086
$ws.Cells.Item(20, 5).Value = 1.040882778643018
$ws.Cells.Item(20, 6).Value = 1.048471311044364
$ws.Cells.Item(20, 9).Value = 1.036144732313139
$ws.Cells.Item(20, 10).Value = 1.039052981763113
$ws.Cells.Item(20, 11).Value = 1.03883861906479
$ws.Cells.Item(20, 12).Value = 1.044239349277612
$ws.Cells.Item(20, 13).Value = 1.051801853073516

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.031873226036362
$ws.Cells.Item(21, 4).Value = 1.034769689527489
$ws.Cells.Item(21, 5).Value = 1.039991956647671
$ws.Cells.Item(21, 6).Value = 1.047401242678897
$ws.Cells.Item(21, 9).Value = 1.035939807291643
$ws.Cells.Item(21, 10).Value = 1.038488874188984
$ws.Cells.Item(21, 11).Value = 1.038354412404526
$ws.Cells.Item(21, 12).Value = 1.043557301372442
$ws.Cells.Item(21, 13).Value = 1.050939460332114

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031268880567684
$ws.Cells.Item(22, 4).Value = 1.034334170547912
$ws.Cells.Item(22, 5).Value = 1.039432937025778
$ws.Cells.Item(22, 6).Value = 1.046729903077206
$ws.Cells.Item(22, 9).Value = 1.035809997678911
$ws.Cells.Item(22, 10).Value = 1.038134160365393
$ws.Cells.Item(22, 11).Value = 1.03804973896979
$ws.Cells.Item(22, 12).Value = 1.043128836964105
$ws.Cells.Item(22, 13).Value = 1.050398024444219

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.031589191483129
$ws.Cells.Item(23, 4).Value = 1.034564997545942
$ws.Cells.Item(23, 5).Value = 1.039729181376466
$ws.Cells.Item(23, 6).Value = 1.047085654405642
$ws.Cells.Item(23, 9).Value = 1.035878903477798
$ws.Cells.Item(23, 10).Value = 1.038322203979167
$ws.Cells.Item(23, 11).Value = 1.038211273862201
$ws.Cells.Item(23, 12).Value = 1.043355938886696
$ws.Cells.Item(23, 13).Value = 1.050684974955321

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.032851637504967
$ws.Cells.Item(24, 4).Value = 1.035474837714719
$ws.Cells.Item(24, 5).Value = 1.04089773149591
$ws.Cells.Item(24, 6).Value = 1.048489275414177
$ws.Cells.Item(24, 9).Value = 1.036148151500061
$ws.Cells.Item(24, 10).Value = 1.039062438385108
$ws.Cells.Item(24, 11).Value = 1.038846732808836
$ws.Cells.Item(24, 12).Value = 1.044250790007637
$ws.Cells.Item(24, 13).Value = 1.051816324388744

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.03431950429127
$ws.Cells.Item(25, 4).Value = 1.036532863532586
$ws.Cells.Item(25, 5).Value = 1.04225839362846
$ws.Cells.Item(25, 6).Value = 1.05012437471542
$ws.Cells.Item(25, 9).Value = 1.036456374819708
$ws.Cells.Item(25, 10).Value = 1.039921237571944
$ws.Cells.Item(25, 11).Value = 1.039583091414209
$ws.Cells.Item(25, 12).Value = 1.045290763391571
$ws.Cells.Item(25, 13).Value = 1.053132556724184

